# Scheduled runner update: refresh cached market-price figures
# (currentAveragePrice / *NQ / *HQ / LevePrice* / LeveProfit*) across
# the per-job leve-profit sheets. Values are plain numeric snapshots
# (no formulas in this workbook), so each touched cell is written
# directly; a couple of profit cells that no longer apply are cleared
# back to blank, and a couple that now apply are populated for the
# first time.
$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 123999.89
$ws.Range("J3").Value = 123999.89
$ws.Range("L3").Value = 123999.89
$ws.Range("N3").Value = -124227.89
$ws.Range("H6").Value = 573.8
$ws.Range("I6").Value = 517.25
$ws.Range("K6").Value = 1551.75
$ws.Range("M6").Value = -1439.75
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -18984
$ws.Range("H102").Value = 123999.89
$ws.Range("J102").Value = 123999.89
$ws.Range("L102").Value = 123999.89
$ws.Range("N102").Value = -130489.89
$ws.Range("H106").Value = 29336312
$ws.Range("I106").Value = 36669196
$ws.Range("K106").Value = 36669196
$ws.Range("M106").Value = -36668565
$ws.Range("H132").Value = 2854.78
$ws.Range("I132").Value = 2603.0889
$ws.Range("K132").Value = 7809.2667
$ws.Range("M132").Value = -5279.2667
$ws.Range("H138").Value = 5606.147
$ws.Range("I138").Value = 4771.3335
$ws.Range("J138").Value = 6265.2104
$ws.Range("K138").Value = 14314.0005
$ws.Range("L138").Value = 18795.6312
$ws.Range("M138").Value = -9174.000499999998
$ws.Range("N138").Value = -29075.6312

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6670.4517
$ws.Range("I32").Value = 2904.8245
$ws.Range("K32").Value = 2904.8245
$ws.Range("M32").Value = -2617.8245
$ws.Range("H61").Value = 8590.637000000001
$ws.Range("I61").Value = 4455.8184
$ws.Range("J61").Value = 16860.273
$ws.Range("K61").Value = 4455.8184
$ws.Range("L61").Value = 16860.273
$ws.Range("M61").Value = -4243.8184
$ws.Range("N61").Value = -17284.273
$ws.Range("H136").Value = 8590.637000000001
$ws.Range("I136").Value = 4455.8184
$ws.Range("J136").Value = 16860.273
$ws.Range("K136").Value = 13367.4552
$ws.Range("L136").Value = 50580.819
$ws.Range("M136").Value = -10817.4552
$ws.Range("N136").Value = -55680.819

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 10000
$ws.Range("J60").Value = 10000
$ws.Range("L60").Value = 10000
$ws.Range("N60").Value = -11198
$ws.Range("H105").Value = 4237.375
$ws.Range("I105").Value = 4291.25
$ws.Range("K105").Value = 4291.25
$ws.Range("M105").Value = -2544.25
$ws.Range("H107").Value = 33334326
$ws.Range("I107").Value = 790.1667
$ws.Range("K107").Value = 790.1667
$ws.Range("M107").Value = 1129.8333
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 267386.1
$ws.Range("I7").Value = 364355.1
$ws.Range("J7").Value = 721.375
$ws.Range("K7").Value = 364355.1
$ws.Range("L7").Value = 721.375
$ws.Range("M7").Value = -364242.1
$ws.Range("N7").Value = -947.375
$ws.Range("H69").Value = 96289.734
$ws.Range("I69").Value = 33680
$ws.Range("J69").Value = 127594.6
$ws.Range("K69").Value = 33680
$ws.Range("L69").Value = 127594.6
$ws.Range("M69").Value = -32931
$ws.Range("N69").Value = -129092.6
$ws.Range("H72").Value = 96289.734
$ws.Range("I72").Value = 33680
$ws.Range("J72").Value = 127594.6
$ws.Range("K72").Value = 101040
$ws.Range("L72").Value = 382783.8
$ws.Range("M72").Value = -97296
$ws.Range("N72").Value = -390271.8
$ws.Range("H93").Value = 131111.6
$ws.Range("I93").Value = 13822.333
$ws.Range("K93").Value = 13822.333
$ws.Range("M93").Value = -11950.333
$ws.Range("H99").Value = 1091408
$ws.Range("I99").Value = 1463887.4
$ws.Range("K99").Value = 1463887.4
$ws.Range("M99").Value = -1462389.4
$ws.Range("H126").Value = 1091408
$ws.Range("I126").Value = 1463887.4
$ws.Range("K126").Value = 4391662.199999999
$ws.Range("M126").Value = -4389192.199999999
$ws.Range("H132").Value = 2754.5667
$ws.Range("I132").Value = 2131.85
$ws.Range("K132").Value = 6395.549999999999
$ws.Range("M132").Value = -3865.549999999999
$ws.Range("H134").Value = 5623.222
$ws.Range("I134").Value = 3126.8333
$ws.Range("K134").Value = 9380.499899999999
$ws.Range("M134").Value = -6845.499899999999

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 99
$ws.Range("I41").Value = 99
$ws.Range("K41").Value = 297
$ws.Range("M41").Value = 41
$ws.Range("H56").Value = 3989.84
$ws.Range("I56").Value = 3989.84
$ws.Range("K56").Value = 3989.84
$ws.Range("M56").Value = -3459.84
$ws.Range("H94").Value = 12011.3
$ws.Range("I94").Value = 4778.5
$ws.Range("J94").Value = 16833.166
$ws.Range("K94").Value = 14335.5
$ws.Range("L94").Value = 50499.49800000001
$ws.Range("M94").Value = -13659.5
$ws.Range("N94").Value = -51851.49800000001
$ws.Range("H109").Value = 1506.75
$ws.Range("I109").Value = 1506.75
$ws.Range("K109").Value = 4520.25
$ws.Range("M109").Value = -3480.25
$ws.Range("H132").Value = 3097.9048
$ws.Range("I132").Value = 3079.4375
$ws.Range("K132").Value = 27714.9375
$ws.Range("M132").Value = -25184.9375

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 157107.14
$ws.Range("J93").Value = 157107.14
$ws.Range("L93").Value = 157107.14
$ws.Range("N93").Value = -160851.14
$ws.Range("H102").Value = 3228.4
$ws.Range("I102").Value = 1574.2667
$ws.Range("K102").Value = 1574.2667
$ws.Range("M102").Value = 47.7333000000001
$ws.Range("H106").Value = 145247
$ws.Range("J106").Value = 145247
$ws.Range("L106").Value = 145247
$ws.Range("N106").Value = -147771
$ws.Range("H122").Value = 3513.3333
$ws.Range("J122").Value = 4885.5
$ws.Range("L122").Value = 14656.5
$ws.Range("N122").Value = -19556.5

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47219.082
$ws.Range("I7").Value = 57919.79
$ws.Range("J7").Value = 6556.4
$ws.Range("K7").Value = 57919.79
$ws.Range("L7").Value = 6556.4
$ws.Range("M7").Value = -57807.79
$ws.Range("N7").Value = -6780.4
$ws.Range("H22").Value = 1471.1428
$ws.Range("I22").Value = 566.6667
$ws.Range("J22").Value = 2149.5
$ws.Range("K22").Value = 566.6667
$ws.Range("L22").Value = 2149.5
$ws.Range("M22").Value = -271.6667
$ws.Range("N22").Value = -2739.5
$ws.Range("H27").Value = 1471.1428
$ws.Range("I27").Value = 566.6667
$ws.Range("J27").Value = 2149.5
$ws.Range("K27").Value = 566.6667
$ws.Range("L27").Value = 2149.5
$ws.Range("M27").Value = -459.6667
$ws.Range("N27").Value = -2363.5
$ws.Range("H35").Value = 26380
$ws.Range("I35").Value = 1725
$ws.Range("K35").Value = 1725
$ws.Range("M35").Value = -1389
$ws.Range("H40").Value = 2772.8857
$ws.Range("I40").Value = 2178.5334
$ws.Range("J40").Value = 6339
$ws.Range("K40").Value = 2178.5334
$ws.Range("L40").Value = 6339
$ws.Range("M40").Value = -2042.5334
$ws.Range("N40").Value = -6611
$ws.Range("H46").Value = 4787.615
$ws.Range("I46").Value = 3992.2856
$ws.Range("J46").Value = 5715.5
$ws.Range("K46").Value = 3992.2856
$ws.Range("L46").Value = 5715.5
$ws.Range("M46").Value = -3804.2856
$ws.Range("N46").Value = -6091.5
$ws.Range("H55").Value = 865.85187
$ws.Range("J55").Value = 1137.5
$ws.Range("L55").Value = 1137.5
$ws.Range("N55").Value = -1483.5
$ws.Range("H58").Value = 11959.4
$ws.Range("J58").Value = 17932.334
$ws.Range("L58").Value = 17932.334
$ws.Range("N58").Value = -18452.334
$ws.Range("H61").Value = 5037.143
$ws.Range("I61").Value = 3228.6365
$ws.Range("J61").Value = 11668.333
$ws.Range("K61").Value = 3228.6365
$ws.Range("L61").Value = 11668.333
$ws.Range("M61").Value = -3026.6365
$ws.Range("N61").Value = -12072.333
$ws.Range("H113").Value = 5037.143
$ws.Range("I113").Value = 3228.6365
$ws.Range("J113").Value = 11668.333
$ws.Range("K113").Value = 3228.6365
$ws.Range("L113").Value = 11668.333
$ws.Range("M113").Value = -1058.6365
$ws.Range("N113").Value = -16008.333
$ws.Range("H126").Value = 47219.082
$ws.Range("I126").Value = 57919.79
$ws.Range("J126").Value = 6556.4
$ws.Range("K126").Value = 173759.37
$ws.Range("L126").Value = 19669.2
$ws.Range("M126").Value = -171289.37
$ws.Range("N126").Value = -24609.2
